$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.427.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.588.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.593.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.02%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.041.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.341.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.586.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.701.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0845"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.848"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "296.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0562"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.009.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
